$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 188, shifting existing rows 188-218 down to 189-219.
$ws.Range("A188:R188").Insert(-4121)

# Copy the formatting (number formats/styles) from the row above (old row 188,
# now shifted to row 189) into the freshly inserted row 188 so the date cell
# keeps its date style.
$ws.Range("A189:R189").Copy()
$ws.Range("A188:R188").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(188, 1).Value = 8
$ws.Cells.Item(188, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(188, 3).Value = "Coquimbo"
$ws.Cells.Item(188, 4).Value = 44504
$ws.Cells.Item(188, 5).Value = 4
$ws.Cells.Item(188, 6).Value = 100114013
$ws.Cells.Item(188, 7).Value = "Zanahoria"
$ws.Cells.Item(188, 8).Value = "Sin especificar"
$ws.Cells.Item(188, 9).Value = "Primera"
$ws.Cells.Item(188, 10).Value = 600
$ws.Cells.Item(188, 11).Value = 6500
$ws.Cells.Item(188, 12).Value = 7000
$ws.Cells.Item(188, 13).Value = 6750
$ws.Cells.Item(188, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(188, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(188, 16).Value = 338
$ws.Cells.Item(188, 17).Value = 20
$ws.Cells.Item(188, 18).Value = "Hortaliza"
